{"js": "// Replace the date line and each \"NNN\u00f7N=\" division-problem cell with its\n// updated value, per the commit's regenerated-worksheet diff. Every old\n// value in this table is unique, so a straightforward search+replace by\n// exact text is safe and unambiguous.\nconst replacements = [\n  [\"2025-08-27 Wednesday\", \"2025-08-28 Thursday\"],\n  [\"148\u00f75=\", \"809\u00f79=\"],\n  [\"201\u00f75=\", \"271\u00f77=\"],\n  [\"429\u00f76=\", \"841\u00f78=\"],\n  [\"918\u00f78=\", \"303\u00f76=\"],\n  [\"338\u00f78=\", \"886\u00f74=\"],\n  [\"327\u00f77=\", \"845\u00f79=\"],\n  [\"526\u00f78=\", \"324\u00f73=\"],\n  [\"916\u00f79=\", \"432\u00f75=\"],\n  [\"207\u00f72=\", \"880\u00f72=\"],\n  [\"203\u00f79=\", \"622\u00f76=\"],\n  [\"592\u00f73=\", \"482\u00f79=\"],\n  [\"724\u00f78=\", \"984\u00f77=\"],\n  [\"141\u00f77=\", \"543\u00f77=\"],\n  [\"654\u00f73=\", \"176\u00f78=\"],\n  [\"759\u00f76=\", \"908\u00f77=\"],\n  [\"602\u00f72=\", \"546\u00f74=\"],\n  [\"762\u00f79=\", \"517\u00f78=\"],\n  [\"651\u00f73=\", \"782\u00f74=\"],\n  [\"587\u00f78=\", \"527\u00f78=\"],\n  [\"233\u00f76=\", \"776\u00f76=\"],\n  [\"593\u00f75=\", \"289\u00f74=\"],\n  [\"478\u00f76=\", \"567\u00f79=\"],\n  [\"697\u00f73=\", \"204\u00f74=\"],\n  [\"461\u00f72=\", \"338\u00f75=\"],\n  [\"945\u00f76=\", \"494\u00f77=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each \"NNN\u00f7N=\" division-problem cell with its\n# updated value, per the commit's regenerated-worksheet diff. Every old\n# value in this table is unique, so Find/Replace on exact text is safe and\n# unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ old = \"2025-08-27 Wednesday\"; new = \"2025-08-28 Thursday\" },\n    @{ old = \"148\u00f75=\"; new = \"809\u00f79=\" },\n    @{ old = \"201\u00f75=\"; new = \"271\u00f77=\" },\n    @{ old = \"429\u00f76=\"; new = \"841\u00f78=\" },\n    @{ old = \"918\u00f78=\"; new = \"303\u00f76=\" },\n    @{ old = \"338\u00f78=\"; new = \"886\u00f74=\" },\n    @{ old = \"327\u00f77=\"; new = \"845\u00f79=\" },\n    @{ old = \"526\u00f78=\"; new = \"324\u00f73=\" },\n    @{ old = \"916\u00f79=\"; new = \"432\u00f75=\" },\n    @{ old = \"207\u00f72=\"; new = \"880\u00f72=\" },\n    @{ old = \"203\u00f79=\"; new = \"622\u00f76=\" },\n    @{ old = \"592\u00f73=\"; new = \"482\u00f79=\" },\n    @{ old = \"724\u00f78=\"; new = \"984\u00f77=\" },\n    @{ old = \"141\u00f77=\"; new = \"543\u00f77=\" },\n    @{ old = \"654\u00f73=\"; new = \"176\u00f78=\" },\n    @{ old = \"759\u00f76=\"; new = \"908\u00f77=\" },\n    @{ old = \"602\u00f72=\"; new = \"546\u00f74=\" },\n    @{ old = \"762\u00f79=\"; new = \"517\u00f78=\" },\n    @{ old = \"651\u00f73=\"; new = \"782\u00f74=\" },\n    @{ old = \"587\u00f78=\"; new = \"527\u00f78=\" },\n    @{ old = \"233\u00f76=\"; new = \"776\u00f76=\" },\n    @{ old = \"593\u00f75=\"; new = \"289\u00f74=\" },\n    @{ old = \"478\u00f76=\"; new = \"567\u00f79=\" },\n    @{ old = \"697\u00f73=\"; new = \"204\u00f74=\" },\n    @{ old = \"461\u00f72=\"; new = \"338\u00f75=\" },\n    @{ old = \"945\u00f76=\"; new = \"494\u00f77=\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair.old, $false, $false, $false, $false, $false, $true, 1, $false, $pair.new, 2) | Out-Null\n}\n"}
